$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H51").Value = 4513.778
$ws.Range("I51").Value = 3441.8333
$ws.Range("J51").Value = 4728.1665
$ws.Range("K51").Value = 3441.8333
$ws.Range("L51").Value = 4728.1665
$ws.Range("M51").Value = -2957.8333
$ws.Range("N51").Value = -5696.1665

$ws = $wb.Worksheets.Item(1)
$ws.Range("H114").Value = 99999
$ws.Range("J114").Value = 99999
$ws.Range("L114").Value = 99999
$ws.Range("N114").Value = -108677

$ws = $wb.Worksheets.Item(1)
$ws.Range("H129").Value = 3906.3333
$ws.Range("I129").Value = 805.25
$ws.Range("J129").Value = 8041.1113
$ws.Range("K129").Value = 2415.75
$ws.Range("L129").Value = 24123.3339
$ws.Range("M129").Value = 2584.25
$ws.Range("N129").Value = -34123.3339

$ws = $wb.Worksheets.Item(1)
$ws.Range("H132").Value = 4408.9785
$ws.Range("I132").Value = 1983.875
$ws.Range("J132").Value = 9582.532999999999
$ws.Range("K132").Value = 5951.625
$ws.Range("L132").Value = 28747.599
$ws.Range("M132").Value = -3421.625
$ws.Range("N132").Value = -33807.599

$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value = 8455.046
$ws.Range("J138").Value = 8096.1
$ws.Range("L138").Value = 24288.3
$ws.Range("N138").Value = -34568.3

$ws = $wb.Worksheets.Item(1)
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 4510.6665
$ws.Range("I32").Value = 4510.6665
$ws.Range("K32").Value = 4510.6665
$ws.Range("M32").Value = -4223.6665

$ws = $wb.Worksheets.Item(2)
$ws.Range("H43").Value = 38656
$ws.Range("J43").Value = 38296.668
$ws.Range("L43").Value = 38296.668
$ws.Range("N43").Value = -38922.668

$ws = $wb.Worksheets.Item(2)
$ws.Range("H110").Value = 9614.333000000001
$ws.Range("I110").Value = 10537.2
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 10537.2
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = -8492.200000000001
$ws.Range("N110").Value = -9090

$ws = $wb.Worksheets.Item(2)
$ws.Range("H132").Value = 2706673.2
$ws.Range("I132").Value = 3955.0908
$ws.Range("J132").Value = 25004098
$ws.Range("K132").Value = 11865.2724
$ws.Range("L132").Value = 75012294
$ws.Range("M132").Value = -9335.2724
$ws.Range("N132").Value = -75017354

$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 5039.769
$ws.Range("I86").Value = 2342.7144
$ws.Range("J86").Value = 8186.3335
$ws.Range("K86").Value = 2342.7144
$ws.Range("L86").Value = 8186.3335
$ws.Range("M86").Value = -1219.7144
$ws.Range("N86").Value = -10432.3335

$ws = $wb.Worksheets.Item(3)
$ws.Range("H89").Value = 5039.769
$ws.Range("I89").Value = 2342.7144
$ws.Range("J89").Value = 8186.3335
$ws.Range("K89").Value = 11713.572
$ws.Range("L89").Value = 40931.6675
$ws.Range("M89").Value = -6097.572
$ws.Range("N89").Value = -52163.6675

$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 4350295.5
$ws.Range("I134").Value = 2550.2
$ws.Range("K134").Value = 7650.599999999999
$ws.Range("M134").Value = -5115.599999999999

$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 839.6667
$ws.Range("I7").Value = 260
$ws.Range("K7").Value = 260
$ws.Range("M7").Value = -147

$ws = $wb.Worksheets.Item(4)
$ws.Range("I132").Value = 2599.48
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 7798.440000000001
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -5268.440000000001
$ws.Range("N132").Value = -17810

$ws = $wb.Worksheets.Item(5)
$ws.Range("H12").Value = 1998.5
$ws.Range("I12").Value = 101.5
$ws.Range("J12").Value = 2757.3
$ws.Range("K12").Value = 304.5
$ws.Range("L12").Value = 8271.900000000001
$ws.Range("M12").Value = -131.5
$ws.Range("N12").Value = -8617.900000000001

$ws = $wb.Worksheets.Item(5)
$ws.Range("H50").Value = 1013.3333
$ws.Range("I50").Value = 1500
$ws.Range("J50").Value = 40
$ws.Range("K50").Value = 4500
$ws.Range("L50").Value = 120
$ws.Range("M50").Value = -4019
$ws.Range("N50").Value = -1082

$ws = $wb.Worksheets.Item(5)
$ws.Range("H53").Value = 1013.3333
$ws.Range("I53").Value = 1500
$ws.Range("J53").Value = 40
$ws.Range("K53").Value = 4500
$ws.Range("L53").Value = 120
$ws.Range("M53").Value = -4019
$ws.Range("N53").Value = -1082

$ws = $wb.Worksheets.Item(5)
$ws.Range("H56").Value = 15317.426
$ws.Range("I56").Value = 15317.426
$ws.Range("K56").Value = 15317.426
$ws.Range("M56").Value = -14787.426

$ws = $wb.Worksheets.Item(5)
$ws.Range("H86").Value = 989.8333
$ws.Range("I86").Value = 399
$ws.Range("J86").Value = 1108
$ws.Range("K86").Value = 1197
$ws.Range("L86").Value = 3324
$ws.Range("M86").Value = -11
$ws.Range("N86").Value = -5696

$ws = $wb.Worksheets.Item(5)
$ws.Range("H89").Value = 989.8333
$ws.Range("I89").Value = 399
$ws.Range("J89").Value = 1108
$ws.Range("K89").Value = 3591
$ws.Range("L89").Value = 9972
$ws.Range("M89").Value = 2337
$ws.Range("N89").Value = -21828

$ws = $wb.Worksheets.Item(5)
$ws.Range("H92").Value = 200
$ws.Range("J92").Value = 200
$ws.Range("L92").Value = 600
$ws.Range("N92").Value = -3096

$ws = $wb.Worksheets.Item(5)
$ws.Range("H110").Value = 33333
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws = $wb.Worksheets.Item(5)
$ws.Range("H111").Value = 18766.5
$ws.Range("I111").Value = 4200
$ws.Range("J111").Value = 33333
$ws.Range("K111").Value = 12600
$ws.Range("L111").Value = 99999
$ws.Range("M111").Value = -9533
$ws.Range("N111").Value = -106133

$ws = $wb.Worksheets.Item(5)
$ws.Range("H112").Value = 18416.5
$ws.Range("I112").Value = 4000
$ws.Range("J112").Value = 20476
$ws.Range("K112").Value = 12000
$ws.Range("L112").Value = 61428
$ws.Range("M112").Value = -10892
$ws.Range("N112").Value = -63644

$ws = $wb.Worksheets.Item(5)
$ws.Range("H132").Value = 2824.3333
$ws.Range("J132").Value = 3248.5
$ws.Range("L132").Value = 29236.5
$ws.Range("N132").Value = -34296.5

$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 2288.1428
$ws.Range("I80").Value = 1913.7778
$ws.Range("K80").Value = 1913.7778
$ws.Range("M80").Value = -915.7778000000001

$ws = $wb.Worksheets.Item(6)
$ws.Range("H83").Value = 2288.1428
$ws.Range("I83").Value = 1913.7778
$ws.Range("K83").Value = 9568.889000000001
$ws.Range("M83").Value = -4576.889000000001

$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 3706203.5
$ws.Range("I113").Value = 2999.5
$ws.Range("K113").Value = 2999.5
$ws.Range("M113").Value = -829.5

$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 2193.6924
$ws.Range("I122").Value = 2231.1
$ws.Range("J122").Value = 2069
$ws.Range("K122").Value = 6693.299999999999
$ws.Range("L122").Value = 6207
$ws.Range("M122").Value = -4243.299999999999
$ws.Range("N122").Value = -11107

$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 12505547
$ws.Range("I132").Value = 5393.1665
$ws.Range("J132").Value = 50006010
$ws.Range("K132").Value = 16179.4995
$ws.Range("L132").Value = 150018030
$ws.Range("M132").Value = -13649.4995
$ws.Range("N132").Value = -150023090

$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 29499.75
$ws.Range("I22").Value = 51500
$ws.Range("J22").Value = 7499.5
$ws.Range("K22").Value = 51500
$ws.Range("L22").Value = 7499.5
$ws.Range("M22").Value = -51205
$ws.Range("N22").Value = -8089.5

$ws = $wb.Worksheets.Item(7)
$ws.Range("H27").Value = 29499.75
$ws.Range("I27").Value = 51500
$ws.Range("J27").Value = 7499.5
$ws.Range("K27").Value = 51500
$ws.Range("L27").Value = 7499.5
$ws.Range("M27").Value = -51393
$ws.Range("N27").Value = -7713.5

$ws = $wb.Worksheets.Item(7)
$ws.Range("H68").Value = 2452430.5
$ws.Range("J68").Value = 2350
$ws.Range("L68").Value = 2350
$ws.Range("N68").Value = -3848

$ws = $wb.Worksheets.Item(7)
$ws.Range("H71").Value = 2452430.5
$ws.Range("J71").Value = 2350
$ws.Range("L71").Value = 11750
$ws.Range("N71").Value = -19238

$ws = $wb.Worksheets.Item(7)
$ws.Range("H132").Value = 5003.8
$ws.Range("I132").Value = 3329.2
$ws.Range("J132").Value = 6678.4
$ws.Range("K132").Value = 9987.599999999999
$ws.Range("L132").Value = 20035.2
$ws.Range("M132").Value = -7457.599999999999
$ws.Range("N132").Value = -25095.2

$ws = $wb.Worksheets.Item(8)
$ws.Range("H42").Value = 49974.5
$ws.Range("J42").Value = 49974.5
$ws.Range("L42").Value = 49974.5
$ws.Range("N42").Value = -50730.5

$ws = $wb.Worksheets.Item(8)
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws = $wb.Worksheets.Item(8)
$ws.Range("H56").Value = 57378.668
$ws.Range("J56").Value = 57378.668
$ws.Range("L56").Value = 57378.668
$ws.Range("N56").Value = -58806.668

$ws = $wb.Worksheets.Item(8)
$ws.Range("H81").Value = 2796.5715
$ws.Range("I81").Value = 1597.6666
$ws.Range("J81").Value = 9990
$ws.Range("K81").Value = 3195.3332
$ws.Range("L81").Value = 19980
$ws.Range("M81").Value = -2134.3332
$ws.Range("N81").Value = -22102

$ws = $wb.Worksheets.Item(8)
$ws.Range("H84").Value = 2796.5715
$ws.Range("I84").Value = 1597.6666
$ws.Range("J84").Value = 9990
$ws.Range("K84").Value = 15976.666
$ws.Range("L84").Value = 99900
$ws.Range("M84").Value = -10672.666
$ws.Range("N84").Value = -110508

$ws = $wb.Worksheets.Item(8)
$ws.Range("H96").Value = 5448.533
$ws.Range("I96").Value = 7781.6665
$ws.Range("J96").Value = 3893.111
$ws.Range("K96").Value = 7781.6665
$ws.Range("L96").Value = 3893.111
$ws.Range("M96").Value = -6408.6665
$ws.Range("N96").Value = -6639.111

$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 10000000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 10000000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 30000000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -30005060

$ws = $wb.Worksheets.Item(8)
$ws.Range("H135").Value = 94989.5
$ws.Range("J135").Value = 94989.5
$ws.Range("L135").Value = 94989.5
$ws.Range("N135").Value = -105129.5
